$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 35 (write image/japanese/english first so new shared strings
# land in this order ahead of the later typo-fix strings)
$ws.Range("E35").Value = "hero_challenges_large_labyrinth.jpg"
$ws.Range("B35").Value = "大迷宮に勇者が挑ソード・ワールドRPGリプレイ集アンマント財宝編2"
$ws.Range("C35").Value = "A Hero Challenges a Large Labyrinth: Sword World Replay Collection 2"
$ws.Range("A35").Value = 1998
$ws.Range("D35").Value = "Fujimi Shobo"

# Fix typo in C34: "... Replace Collection 1" -> "... Replay Collection 1"
$ws.Range("C34").Value = "Heroes Gather at the Treasure Map: Sword World Replay Collection 1"

# Fix typo in C33: "Leave it to the Adventurer! ..." -> "Leave It To the Adventurer: ..."
$ws.Range("C33").Value = "Leave It To the Adventurer: Sword World Replay Collection 2"

# Move selection below the newly-added row, matching the saved workbook state
$ws.Range("A36").Select()
